# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columns D ("provincia") and K ("explotaciones-cuya-gestion-se-lleva-por-otra-persona")
# move from being curated as "dim" (dimension) columns to "medida" (measure) columns:
#   - Row 2 (semantic type): sdmx-dimension:refArea -> iaest-measure:provincia
#                             iaest-dimension:...otra-persona -> iaest-measure:...otra-persona
#   - Row 3 (dim/medida flag): dim -> medida
#   - Row 4 (datatype): URI-Provincia / skos:Concept -> xsd:int
#   - Row 5 (mapping file): the mapping xlsx reference is removed (cell cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("provincia") - row 2: refArea -> measure:provincia
$ws.Range("D2").Value = "iaest-measure:provincia"
# Column D - row 3: dim -> medida
$ws.Range("D3").Value = "medida"
# Column D - row 4: URI-Provincia -> xsd:int
$ws.Range("D4").Value = "xsd:int"

# Column K ("explotaciones-cuya-gestion-se-lleva-por-otra-persona") - row 2: dimension -> measure
$ws.Range("K2").Value = "iaest-measure:explotaciones-cuya-gestion-se-lleva-por-otra-persona"
# Column K - row 3: dim -> medida
$ws.Range("K3").Value = "medida"
# Column K - row 4: skos:Concept -> xsd:int
$ws.Range("K4").Value = "xsd:int"
# Column K - row 5: mapping file reference removed (no longer a dimension, no mapping needed)
$ws.Range("K5").Clear()
